$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (columns B..Q) that apply identically to every data row (2..26)
$values = @(
    [double]"0.9999993052626588",    # B r2
    [double]"0.9990633197468317",    # C r2_sup
    [double]"0.9999987880483845",    # D r2_test
    [double]"0.9999999999999456",    # E r2_val
    [double]"0.9999994371058077",    # F r2_vt
    [double]"6.485068126221803e-07", # G mse
    [double]"0.0008743499008213532", # H mse_sup
    [double]"1.371014771879303e-06", # I mse_test
    [double]"4.713938294857159e-14", # J mse_val
    [double]"6.855074095093429e-07", # K mse_vt
    [double]"4.523788428736971e-05", # L mape
    [double]"0.0008052992068928047", # M rmse
    [double]"1.00000185263291",      # N r2_adj
    [double]"0.000839582461329401",  # O rsd
    [double]"94.49718665609024",     # P aic
    [double]"134.7200888767409"      # Q bic
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
